# "Generate Report for Handoff"
#
# Refreshes the localization-status report after a new handoff round:
#   - Overview sheet: zh-cn / de-de status + timestamp columns move on
#     ("Handed back: in sync with en-US" -> "In Translation" / "Ready for handoff")
#   - zh-cn / de-de detail sheets: Status + Latest Handoff Datetime columns
#     updated, and the new "Ready for handoff" file now carries a
#     "version not latest" Error Detail message.
#   - A few columns are narrowed / widened to fit the refreshed content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
# Row 2 -> 1169e9a6-...md : now "In Translation" as of 2016-09-06 15:50:38
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("G2").Value = "2016-09-06 15:50:38"

# Row 3 -> 32c776b7-...md : now "Ready for handoff" as of 2016-09-06 15:50:38
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 15:50:38"

# Column widths: E/F got narrower (long placeholder text no longer needed)
# (target stored width is 17.2159881591797; compensate for the engine's
# internal +5/6 padding so the round-tripped width lands as close as
# possible to that target)
$overview.Columns.Item(5).ColumnWidth = 16.3826548258464
$overview.Columns.Item(6).ColumnWidth = 16.3826548258464

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 15:50:13"
$zhcn.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15f987df2a43a29062e5ad22e1dce8d3230113d/e2e/1169e9a6-ecc9-4d69-a72d-10aa163b9c7a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c522a1815b23cf4642a098b5c456662397e5fbe6/e2e/1169e9a6-ecc9-4d69-a72d-10aa163b9c7a.md."

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-06 15:50:13"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15f987df2a43a29062e5ad22e1dce8d3230113d/e2e/32c776b7-2827-48d6-965d-1a32db1ec7be.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c522a1815b23cf4642a098b5c456662397e5fbe6/e2e/32c776b7-2827-48d6-965d-1a32db1ec7be.md."

# Column widths: C narrower, P (Error Detail) widened to fit the new message
$zhcn.Columns.Item(3).ColumnWidth = 16.3826548258464
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 15:50:38"
$dede.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15f987df2a43a29062e5ad22e1dce8d3230113d/e2e/1169e9a6-ecc9-4d69-a72d-10aa163b9c7a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c522a1815b23cf4642a098b5c456662397e5fbe6/e2e/1169e9a6-ecc9-4d69-a72d-10aa163b9c7a.md."

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-06 15:50:38"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e15f987df2a43a29062e5ad22e1dce8d3230113d/e2e/32c776b7-2827-48d6-965d-1a32db1ec7be.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c522a1815b23cf4642a098b5c456662397e5fbe6/e2e/32c776b7-2827-48d6-965d-1a32db1ec7be.md."

# Column widths: C narrower, P (Error Detail) widened to fit the new message
$dede.Columns.Item(3).ColumnWidth = 16.3826548258464
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
